# Sync attendance_reports: normalize "Recorded By" (column G) entries so that
# when the value starts with "System", that token is swapped with the last
# token in the comma-separated list (text case preserved as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }
    if ($parts[0] -ne "System") { continue }

    $last = $parts.Count - 1
    $tmp = $parts[0]
    $parts[0] = $parts[$last]
    $parts[$last] = $tmp

    $newVal = [string]::Join(", ", $parts)
    $cell.Value = $newVal
}
